$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 501
$ws.Range("I43").Value = 501
$ws.Range("K43").Value = 501
$ws.Range("M43").Value = -432
$ws.Range("H124").Value = 43780
$ws.Range("J124").Value = 43780
$ws.Range("L124").Value = 43780
$ws.Range("N124").Value = -53600
$ws.Range("H125").Value = 3641.5
$ws.Range("I125").Value = 5126.4
$ws.Range("J125").Value = 1166.6666
$ws.Range("K125").Value = 46137.6
$ws.Range("L125").Value = 10499.9994
$ws.Range("M125").Value = -43677.6
$ws.Range("N125").Value = -15419.9994
$ws.Range("H132").Value = 3585682.2
$ws.Range("I132").Value = 1144.3798
$ws.Range("J132").Value = 23812718
$ws.Range("K132").Value = 3433.1394
$ws.Range("L132").Value = 71438154
$ws.Range("M132").Value = -903.1394
$ws.Range("N132").Value = -71443214
$ws.Range("H138").Value = 3417.2373
$ws.Range("I138").Value = 1730
$ws.Range("J138").Value = 4420.4595
$ws.Range("K138").Value = 5190
$ws.Range("L138").Value = 13261.3785
$ws.Range("M138").Value = -50
$ws.Range("N138").Value = -23541.3785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13859.539
$ws.Range("I32").Value = 9912.044
$ws.Range("J32").Value = 27478.4
$ws.Range("K32").Value = 9912.044
$ws.Range("L32").Value = 27478.4
$ws.Range("M32").Value = -9625.044
$ws.Range("N32").Value = -28052.4
$ws.Range("H61").Value = 235778.84
$ws.Range("I61").Value = 3087.5312
$ws.Range("J61").Value = 912699
$ws.Range("K61").Value = 3087.5312
$ws.Range("L61").Value = 912699
$ws.Range("M61").Value = -2875.5312
$ws.Range("N61").Value = -913123
$ws.Range("H74").Value = 7463835
$ws.Range("I74").Value = 840.1905
$ws.Range("J74").Value = 20001666
$ws.Range("K74").Value = 840.1905
$ws.Range("L74").Value = 20001666
$ws.Range("M74").Value = 33.80949999999996
$ws.Range("N74").Value = -20003414
$ws.Range("H77").Value = 7463835
$ws.Range("I77").Value = 840.1905
$ws.Range("J77").Value = 20001666
$ws.Range("K77").Value = 4200.9525
$ws.Range("L77").Value = 100008330
$ws.Range("M77").Value = 167.0474999999997
$ws.Range("N77").Value = -100017066
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 1926964.1
$ws.Range("I132").Value = 2016.8667
$ws.Range("J132").Value = 4551892.5
$ws.Range("K132").Value = 6050.6001
$ws.Range("L132").Value = 13655677.5
$ws.Range("M132").Value = -3520.6001
$ws.Range("N132").Value = -13660737.5
$ws.Range("H136").Value = 235778.84
$ws.Range("I136").Value = 3087.5312
$ws.Range("J136").Value = 912699
$ws.Range("K136").Value = 9262.5936
$ws.Range("L136").Value = 2738097
$ws.Range("M136").Value = -6712.5936
$ws.Range("N136").Value = -2743197

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 25832.887
$ws.Range("I134").Value = 3187.973
$ws.Range("J134").Value = 145527.42
$ws.Range("K134").Value = 9563.919
$ws.Range("L134").Value = 436582.26
$ws.Range("M134").Value = -7028.919
$ws.Range("N134").Value = -441652.26
$ws.Range("H137").Value = 69604
$ws.Range("I137").Value = 110000
$ws.Range("J137").Value = 59505
$ws.Range("K137").Value = 110000
$ws.Range("L137").Value = 59505
$ws.Range("M137").Value = -104900
$ws.Range("N137").Value = -69705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 1100
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H12").Value = 1417.1666
$ws.Range("I12").Value = 1417.1666
$ws.Range("K12").Value = 1417.1666
$ws.Range("M12").Value = -1247.1666
$ws.Range("H13").Value = 28402
$ws.Range("J13").Value = 31002.223
$ws.Range("L13").Value = 31002.223
$ws.Range("N13").Value = -31280.223
$ws.Range("H31").Value = 15635170
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 15635170
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 15635170
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -15635760
$ws.Range("H34").Value = 15635170
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 15635170
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 15635170
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -15635574
$ws.Range("H58").Value = 5197473
$ws.Range("I58").Value = 6668002.5
$ws.Range("J58").Value = 1112669.5
$ws.Range("K58").Value = 6668002.5
$ws.Range("L58").Value = 1112669.5
$ws.Range("M58").Value = -6667799.5
$ws.Range("N58").Value = -1113075.5
$ws.Range("H86").Value = 1666.8422
$ws.Range("I86").Value = 1699
$ws.Range("J86").Value = 1631.1111
$ws.Range("K86").Value = 1699
$ws.Range("L86").Value = 1631.1111
$ws.Range("M86").Value = -576
$ws.Range("N86").Value = -3877.1111
$ws.Range("H89").Value = 1666.8422
$ws.Range("I89").Value = 1699
$ws.Range("J89").Value = 1631.1111
$ws.Range("K89").Value = 8495
$ws.Range("L89").Value = 8155.5555
$ws.Range("M89").Value = -2879
$ws.Range("N89").Value = -19387.5555
$ws.Range("H134").Value = 7631856.5
$ws.Range("I134").Value = 8549459
$ws.Range("J134").Value = 1667437.5
$ws.Range("K134").Value = 25648377
$ws.Range("L134").Value = 5002312.5
$ws.Range("M134").Value = -25645842
$ws.Range("N134").Value = -5007382.5
$ws.Range("H136").Value = 5197473
$ws.Range("I136").Value = 6668002.5
$ws.Range("J136").Value = 1112669.5
$ws.Range("K136").Value = 20004007.5
$ws.Range("L136").Value = 3338008.5
$ws.Range("M136").Value = -20001457.5
$ws.Range("N136").Value = -3343108.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3368832.5
$ws.Range("I5").Value = 640.8571
$ws.Range("J5").Value = 6996115.5
$ws.Range("K5").Value = 1922.5713
$ws.Range("L5").Value = 20988346.5
$ws.Range("M5").Value = -1810.5713
$ws.Range("N5").Value = -20988570.5
$ws.Range("H129").Value = 1513.75
$ws.Range("I129").Value = 1433.75
$ws.Range("J129").Value = 1593.75
$ws.Range("K129").Value = 4301.25
$ws.Range("L129").Value = 4781.25
$ws.Range("M129").Value = 698.75
$ws.Range("N129").Value = -14781.25
$ws.Range("H131").Value = 1755530.8
$ws.Range("J131").Value = 1396.3077
$ws.Range("L131").Value = 4188.9231
$ws.Range("N131").Value = -14268.9231
$ws.Range("H135").Value = 3368832.5
$ws.Range("I135").Value = 640.8571
$ws.Range("J135").Value = 6996115.5
$ws.Range("K135").Value = 5767.7139
$ws.Range("L135").Value = 62965039.5
$ws.Range("M135").Value = -3232.7139
$ws.Range("N135").Value = -62970109.5
$ws.Range("H140").Value = 5496.4
$ws.Range("I140").Value = 5496.4
$ws.Range("K140").Value = 16489.2
$ws.Range("M140").Value = -11309.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H22").Value = 254
$ws.Range("I22").Value = 254
$ws.Range("K22").Value = 254
$ws.Range("M22").Value = 275
$ws.Range("H132").Value = 16680989
$ws.Range("I132").Value = 20844786
$ws.Range("J132").Value = 25799.5
$ws.Range("K132").Value = 62534358
$ws.Range("L132").Value = 77398.5
$ws.Range("M132").Value = -62531828
$ws.Range("N132").Value = -82458.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 66669304
$ws.Range("I40").Value = 83335890
$ws.Range("J40").Value = 2968.3333
$ws.Range("K40").Value = 83335890
$ws.Range("L40").Value = 2968.3333
$ws.Range("M40").Value = -83335754
$ws.Range("N40").Value = -3240.3333
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42246
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -131232
$ws.Range("H132").Value = 11630684
$ws.Range("I132").Value = 12197839
$ws.Range("J132").Value = 4002.5
$ws.Range("K132").Value = 36593517
$ws.Range("L132").Value = 12007.5
$ws.Range("M132").Value = -36590987
$ws.Range("N132").Value = -17067.5
$ws.Range("H135").Value = 175431.6
$ws.Range("J135").Value = 175431.6
$ws.Range("L135").Value = 175431.6
$ws.Range("N135").Value = -185571.6
$ws.Range("H136").Value = 5419.5415
$ws.Range("I136").Value = 2000.0286
$ws.Range("J136").Value = 14625.923
$ws.Range("K136").Value = 6000.085800000001
$ws.Range("L136").Value = 43877.769
$ws.Range("M136").Value = -3450.085800000001
$ws.Range("N136").Value = -48977.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 755.3214
$ws.Range("I126").Value = 680.4091
$ws.Range("J126").Value = 1030
$ws.Range("K126").Value = 2041.2273
$ws.Range("L126").Value = 3090
$ws.Range("M126").Value = 428.7727
$ws.Range("N126").Value = -8030
$ws.Range("H132").Value = 1383.9354
$ws.Range("I132").Value = 596.16
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 1788.48
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = 741.52
$ws.Range("N132").Value = -19059.0005
$ws.Range("H136").Value = 5559190
$ws.Range("I136").Value = 4458.533
$ws.Range("J136").Value = 11113921
$ws.Range("K136").Value = 13375.599
$ws.Range("L136").Value = 33341763
$ws.Range("M136").Value = -10825.599
$ws.Range("N136").Value = -33346863
